$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.242.94"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "2.171.88"
$ws.Range("E3").Value = "  -1.73%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'236.70"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("D7").Value = "'69.45"
$ws.Range("E7").Value = "  -4.52%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.575"
$ws.Range("E9").Value = "  -4.11%  "
$ws.Range("D10").Value = "'39.40"
$ws.Range("E10").Value = "  -6.09%  "
$ws.Range("D11").Value = "'0.0918"
$ws.Range("E11").Value = "  -3.07%  "
$ws.Range("E12").Value = "  -5.77%  "
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("E14").Value = "  -4.48%  "
$ws.Range("D15").Value = "2.496.05"
$ws.Range("E15").Value = "  -1.78%  "
$ws.Range("D16").Value = "'13.94"
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("D17").Value = "2.174.03"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").Value = "'0.794"
$ws.Range("E18").Value = "  -4.41%  "
$ws.Range("D19").Value = "41.029.49"
$ws.Range("E19").Value = "  -1.96%  "
$ws.Range("D20").Value = "0.0₂01000"
$ws.Range("E20").Value = "  -6.91%  "
$ws.Range("D21").Value = "'70.54"
$ws.Range("E21").Value = "  -3.14%  "
$ws.Range("D22").Value = "'5.88"
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").Value = "'225.49"
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("D24").Value = "'9.35"
$ws.Range("E24").Value = "  -7.78%  "
$ws.Range("E25").Value = "  -7.62%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'10.72"
$ws.Range("E27").Value = "  -7.23%  "
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  -2.55%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'167.92"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("D32").Value = "'19.88"
$ws.Range("E32").Value = "  -3.00%  "
$ws.Range("D33").Value = "'30.57"
$ws.Range("E33").Value = "  +7.10%  "
$ws.Range("D34").Value = "'0.0761"
$ws.Range("E34").Value = "  -2.97%  "
$ws.Range("E35").Value = "  -9.74%  "
$ws.Range("E36").Value = "  -3.08%  "
$ws.Range("E37").Value = "  -7.99%  "
$ws.Range("D38").Value = "'4.05"
$ws.Range("E38").Value = "  -3.72%  "
$ws.Range("D39").Value = "'0.0282"
$ws.Range("E39").Value = "  -5.25%  "
$ws.Range("D40").Value = "'2.06"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").Value = "'5.40"
$ws.Range("E41").Value = "  -3.46%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'11.51"
$ws.Range("E42").Value = "  -12.94%  "
$ws.Range("D43").Value = "'58.94"
$ws.Range("E43").Value = "  -9.06%  "
$ws.Range("E44").Value = "  -3.53%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'8.24"
$ws.Range("E45").Value = "  -4.61%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.0967"
$ws.Range("E46").Value = "  -3.34%  "
$ws.Range("D47").Value = "'97.37"
$ws.Range("E47").Value = "  -5.78%  "
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("D49").Value = "'1.12"
$ws.Range("E49").Value = "  -3.45%  "
$ws.Range("E50").Value = "  -6.97%  "
$ws.Range("E51").Value = "  -2.56%  "
